# Apply edits to Junction_Flooding_314.xlsx-style workbook:
#  - Round all numeric values in row 5 (columns C..Q, T..AH) to 2 decimal places
#  - Delete row 6 entirely (shifting nothing below, as it is the last row)
#  - The dimension will recompute automatically once row 6 is removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that changed precision in row 5 (B, F, R, S unchanged - already <=2 decimals)
$cols = @("C","D","E","G","H","I","J","K","L","M","N","O","P","Q","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")

foreach ($col in $cols) {
    $cell = $ws.Range($col + "5")
    $val = $cell.Value2
    $rounded = [Math]::Round($val, 2)
    $cell.Value = $rounded
}

# Delete row 6 (the last data row) entirely
$ws.Rows.Item(6).Delete()
